# Update the progress/description notes for the Customer and Bills manager tasks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = "Class Creation done (see src) and manager almost finish (not sure about some functions and it misses the connection with the database) if you can check my code to tell me my mistake…"
$ws.Range("E8").Value = "Class Creation done (see src) and manager in progress (group with Customer and upload by Yamid later)"

# Match the saved active selection recorded in the workbook (cell F17).
$ws.Range("F17").Select()
